# Weekly update: insert two new daily-price rows at the top of the
# Brocoli / Femacal de La Calera data block (sheet rows 1065-1066),
# pushing the existing historical rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 1065; everything
# that used to live at row 1065 onward moves down to row 1067 onward.
$ws.Rows("1065:1066").Insert()

# --- New row 1065: "Primera" quality entry for the new date ---
$ws.Cells.Item(1065, 1).Value = 3
$ws.Cells.Item(1065, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1065, 3).Value = "Coquimbo"
$ws.Cells.Item(1065, 4).Value = 45075
$ws.Cells.Item(1065, 5).Value = 5
$ws.Cells.Item(1065, 6).Value = 100112023
$ws.Cells.Item(1065, 7).Value = "Brócoli"
$ws.Cells.Item(1065, 8).Value = "Sin especificar"
$ws.Cells.Item(1065, 9).Value = "Primera"
$ws.Cells.Item(1065, 10).Value = 3750
$ws.Cells.Item(1065, 11).Value = 900
$ws.Cells.Item(1065, 12).Value = 950
$ws.Cells.Item(1065, 13).Value = 924
$ws.Cells.Item(1065, 14).Value = "`$/unidad"
$ws.Cells.Item(1065, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1065, 16).Value = 924
$ws.Cells.Item(1065, 17).Value = 1
$ws.Cells.Item(1065, 18).Value = "Hortaliza"

# --- New row 1066: "Segunda" quality entry for the same new date ---
$ws.Cells.Item(1066, 1).Value = 3
$ws.Cells.Item(1066, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1066, 3).Value = "Coquimbo"
$ws.Cells.Item(1066, 4).Value = 45075
$ws.Cells.Item(1066, 5).Value = 5
$ws.Cells.Item(1066, 6).Value = 100112023
$ws.Cells.Item(1066, 7).Value = "Brócoli"
$ws.Cells.Item(1066, 8).Value = "Sin especificar"
$ws.Cells.Item(1066, 9).Value = "Segunda"
$ws.Cells.Item(1066, 10).Value = 1800
$ws.Cells.Item(1066, 11).Value = 800
$ws.Cells.Item(1066, 12).Value = 800
$ws.Cells.Item(1066, 13).Value = 800
$ws.Cells.Item(1066, 14).Value = "`$/unidad"
$ws.Cells.Item(1066, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1066, 16).Value = 800
$ws.Cells.Item(1066, 17).Value = 1
$ws.Cells.Item(1066, 18).Value = "Hortaliza"
